$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 updates
$ws.Range("G2").Value = 3.6
$ws.Range("I2").Value = 2.45
$ws.Range("W2").Value = 7
$ws.Range("AH2").Value = 9.5
$ws.Range("AQ2").Value = 101
$ws.Range("AZ2").Value = 51

# Row 3 updates
$ws.Range("S3").Value = 1.33

# Row 4 updates
$ws.Range("S4").Value = 1.44
$ws.Range("T4").Value = 2.63
